# Generate Report for Handoff
# Inserts a new status row for "15f878de-062b-406a-814c-cb86cd71896e" between the
# existing "9cd5f78b-9f94-4b2a-be44-1835a27397f7" row and the
# "89db2ebc-63e4-416d-a9bf-b3c90d5f0183" row on every sheet (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$newUuid = "15f878de-062b-406a-814c-cb86cd71896e"
$newHash = "e4c4f2a78678fd9e9fad1749b14a5d9c15ea274c"
$insertRow = 7

# ---------------------------------------------------------------------------
# Helper: capture every hyperlink on a sheet (row, column, address, display),
# then remove them all so we can rebuild them after the row insert (the
# engine does not automatically re-home hyperlinks when rows shift).
# ---------------------------------------------------------------------------
function Get-LinkInfo($ws) {
    $links = @()
    foreach ($hh in $ws.Hyperlinks) {
        $links += @{
            row = $hh.Range.Row
            col = $hh.Range.Column
            addr = $hh.Address
            text = $hh.TextToDisplay
        }
    }
    return $links
}

# =============================================================================
# Sheet "Overview"
# =============================================================================
$ws1 = $wb.Worksheets.Item("Overview")
$links1 = Get-LinkInfo $ws1
$ws1.Hyperlinks.Delete()

$ws1.Rows.Item($insertRow).Insert()

$ws1.Range("A7").Value = "$newUuid.md"
$ws1.Range("B7").Value = "Ready for handoff"
$ws1.Range("C7").Value = "Ready for handoff"
$ws1.Range("D7").Value = "2016-03-24 11:35:19"

foreach ($lk in $links1) {
    $r = $lk.row
    if ($r -ge $insertRow) { $r = $r + 1 }
    $anchor = $ws1.Cells.Item($r, $lk.col)
    $ws1.Hyperlinks.Add($anchor, $lk.addr, "", "", $lk.text) | Out-Null
}
$ws1.Hyperlinks.Add($ws1.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/$newHash/e2e/$newUuid.md", "", "", "$newUuid.md") | Out-Null

# =============================================================================
# Sheet "zh-cn"
# =============================================================================
$ws2 = $wb.Worksheets.Item("zh-cn")
$links2 = Get-LinkInfo $ws2
$ws2.Hyperlinks.Delete()

$ws2.Rows.Item($insertRow).Insert()

$ws2.Range("A7").Value = "$newUuid.md"
$ws2.Range("B7").Value = ".md"
$ws2.Range("C7").Value = "Ready for handoff"
$ws2.Range("D7").Value = "$newUuid.$newHash.zh-cn.xlf"
$ws2.Range("E7").Value = "2016-03-24 11:35:14"
$ws2.Range("H7").Value = "0001-01-01 00:00:00"
$ws2.Range("J7").Value = "Include"

foreach ($lk in $links2) {
    $r = $lk.row
    if ($r -ge $insertRow) { $r = $r + 1 }
    $anchor = $ws2.Cells.Item($r, $lk.col)
    $ws2.Hyperlinks.Add($anchor, $lk.addr, "", "", $lk.text) | Out-Null
}
$ws2.Hyperlinks.Add($ws2.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/$newHash/e2e/$newUuid.md", "", "", "$newUuid.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$newHash/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newUuid.$newHash.zh-cn.xlf", "", "", "$newUuid.$newHash.zh-cn.xlf") | Out-Null

# =============================================================================
# Sheet "de-de"
# =============================================================================
$ws3 = $wb.Worksheets.Item("de-de")
$links3 = Get-LinkInfo $ws3
$ws3.Hyperlinks.Delete()

$ws3.Rows.Item($insertRow).Insert()

$ws3.Range("A7").Value = "$newUuid.md"
$ws3.Range("B7").Value = ".md"
$ws3.Range("C7").Value = "Ready for handoff"
$ws3.Range("D7").Value = "$newUuid.$newHash.de-de.xlf"
$ws3.Range("E7").Value = "2016-03-24 11:35:19"
$ws3.Range("H7").Value = "0001-01-01 00:00:00"
$ws3.Range("J7").Value = "Include"

foreach ($lk in $links3) {
    $r = $lk.row
    if ($r -ge $insertRow) { $r = $r + 1 }
    $anchor = $ws3.Cells.Item($r, $lk.col)
    $ws3.Hyperlinks.Add($anchor, $lk.addr, "", "", $lk.text) | Out-Null
}
$ws3.Hyperlinks.Add($ws3.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/$newHash/e2e/$newUuid.md", "", "", "$newUuid.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$newHash/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newUuid.$newHash.de-de.xlf", "", "", "$newUuid.$newHash.de-de.xlf") | Out-Null

Write-Host "Done: inserted handoff row for $newUuid on all sheets."
